# Adds support for IP ranges and port ranges:
#  - NetworkAddresses: two new rows (test29 / test30) holding IP ranges
#  - Services: one new row (test_range) holding a tcp port range
# Also widens the NetworkAddresses "Ipaddress" column so the longer
# range strings are readable, and leaves the UI selection/active-sheet
# state the way the author left it (Services ends up the active tab).

$wb = $excel.ActiveWorkbook

$wsFirewall = $wb.Worksheets.Item("FirewallPolicies")
$wsAddresses = $wb.Worksheets.Item("NetworkAddresses")
$wsObjects = $wb.Worksheets.Item("NetworkObjects")
$wsServices = $wb.Worksheets.Item("Services")
$wsGroups = $wb.Worksheets.Item("ServiceGroups")

# --- NetworkAddresses: two new IP-range entries ------------------------
$wsAddresses.Range("A30").Value = "test29"
$wsAddresses.Range("A31").Value = "test30"
$wsAddresses.Range("B30").Value = "1.1.1.1-10.10.10.10"
$wsAddresses.Range("B31").Value = "20.20.20.20-30.30.30.30"

# Widen column B ("Ipaddress") so the new range values fit.
$wsAddresses.Columns.Item(2).ColumnWidth = 24.26953125

# --- Services: new tcp port-range entry --------------------------------
$wsServices.Range("A5").Value = "test_range"
$wsServices.Range("B5").Value = "tcp"
$wsServices.Range("C5").Value = "1000-2000"

# --- Restore the cursor/selection on each sheet -------------------------
# (selections are applied in the order the author last touched the
# sheets, so the final Activate leaves "Services" as the active tab)
$wsFirewall.Range("D7").Select()
$wsAddresses.Range("D24").Select()
$wsGroups.Range("F12").Select()

$wsServices.Activate()
$wsServices.Range("A6").Select()
